$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0157781932502985
$ws.Range("C2").Value = 0.0073637580499053
$ws.Range("D2").Value = 0.0047916327603161335
$ws.Range("E2").Value = 0.0043869344517588615
$ws.Range("F2").Value = 0.0001434998121112585
$ws.Range("G2").Value = 0.0009819269180297852
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.12730003893375397
$ws.Range("K2").Value = 1.418911099433899
